# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the Thu Jul 20 07:34:57 UTC 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.211.66"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = "'1.912.54"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'0.8215"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('D6').Value = "'243.78"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.28%  '
$ws.Range('D7').Value = "'1.000"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.3258"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.72%  '
$ws.Range('D9').Value = "'26.85"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.73%  '
$ws.Range('D10').Value = "'0.07063"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('D11').Value = "'0.08101"
$ws.Range('D11').ClearFormats()
$ws.Range('E12').Value = '  +3.34%  '
$ws.Range('D13').Value = "'1.919.26"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = "'5.296"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.06%  '
$ws.Range('D15').Value = "'93.46"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = "'30.211.08"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').Value = "'5.928"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.29%  '
$ws.Range('D19').Value = "'246.74"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = "'0.000007802"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = "'2.165.76"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = "'1.001"
$ws.Range('D22').ClearFormats()
$ws.Range('D23').Value = "'1.000"
$ws.Range('D23').ClearFormats()
$ws.Range('D24').Value = "'7.107"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.71%  '
$ws.Range('D25').Value = "'0.1678"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +20.48%  '
$ws.Range('D26').Value = "'9.340"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').Value = "'167.29"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('D28').Value = "'19.01"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('E29').Value = '  +3.23%  '
$ws.Range('D30').Value = "'1.374"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = "'1.529"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('D32').Value = "'0.05865"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.80%  '
$ws.Range('D33').Value = "'4.315"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('D34').Value = "'4.110"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D35').Value = "'1.277"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.00%  '
$ws.Range('D36').Value = "'0.7377"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('D38').Value = "'0.01926"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('D39').Value = "'2.800"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').Value = "'0.4470"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.44%  '
$ws.Range('D41').Value = "'73.48"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.12%  '
$ws.Range('D42').Value = "'5.979"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('D43').Value = "'0.8539"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.21%  '
$ws.Range('D44').Value = "'1.918"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').Value = "'0.9999"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = "'102.92"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('D47').Value = "'7.602"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = "'9.881"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('D49').Value = "'1.009.35"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('D50').Value = "'2.064.12"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = "'1.559"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.87%  '
